# 文档 by qiyu 2023-03-01
# Append five new "Title and Content" slides to the end of the deck.

$p = $ppt.ActivePresentation

$newSlides = @(
    @{ En = "Cypher"; Zh = "执行计划" },
    @{ En = "Neo4j";  Zh = "免邻接索引" },
    @{ En = "Neo4j";  Zh = "底层存储结构" },
    @{ En = "Neo4j";  Zh = "索引，约束" },
    @{ En = "Neo4j";  Zh = "集群" }
)

foreach ($item in $newSlides) {
    $idx = $p.Slides.Count + 1
    $s = $p.Slides.Add($idx, 2)

    $titleShape = $s.Shapes.Item(1)
    $titleShape.Name = "标题 1"
    $titleRange = $titleShape.TextFrame.TextRange
    $titleRange.Text = $item.En
    $titleRange.InsertAfter($item.Zh) | Out-Null

    $bodyShape = $s.Shapes.Item(2)
    $bodyShape.Name = "内容占位符 2"
}
